$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell address -> new text value, taken from the refreshed symbol-list scrape.
# Values are written verbatim as text (matching the workbooks existing inlineStr
# columns), so things like trailing zeros ("0.03440") and percent suffixes ("5.68%")
# must survive exactly rather than being renormalised as numbers.
$updates = [ordered]@{
    'D2' = '329.03'
    'E2' = '5.68%'
    'D3' = '40.04'
    'E3' = '6.45%'
    'D4' = '5.274'
    'E4' = '2.14%'
    'D5' = '0.08113'
    'E5' = '2.34%'
    'B6' = 'GateToken'
    'C6' = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
    'D6' = '4.529'
    'E6' = '2.26%'
    'B7' = 'KuCoinToken'
    'C7' = 'https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs'
    'D7' = '8.640'
    'E7' = '4.32%'
    'B8' = 'FTXToken'
    'C8' = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
    'D8' = '1.914'
    'E8' = '-0.53%'
    'B9' = 'BTSEToken'
    'C9' = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
    'D9' = '2.958'
    'E9' = '-1.41%'
    'B10' = 'MXToken'
    'C10' = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
    'D10' = '0.9346'
    'E10' = '0.46%'
    'B11' = 'LiechtensteinCryptoassetsExchange'
    'C11' = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
    'D11' = '0.1344'
    'E11' = '22.40%'
    'B12' = 'WazirX'
    'C12' = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
    'D12' = '0.1956'
    'E12' = '1.90%'
    'B13' = 'MandalaExchangeToken'
    'C13' = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
    'D13' = '0.09179'
    'E13' = '1.14%'
    'B14' = 'BitrueCoin'
    'C14' = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
    'D14' = '0.03440'
    'E14' = '4.24%'
    'B15' = 'BitMartToken'
    'C15' = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
    'D15' = '0.09558'
    'E15' = '-0.33%'
    'B16' = 'BitForexToken'
    'C16' = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
    'D16' = '0.001390'
    'E16' = '1.11%'
    'B17' = 'TigerCash'
    'C17' = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
    'D17' = '0.006013'
    'E17' = '3.01%'
    'B18' = 'LEO'
    'C18' = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
    'D18' = '3.360'
    'E18' = '-6.51%'
    'D19' = '0.3526'
    'E19' = '3.42%'
    'D20' = '7.276'
    'E20' = '21.74%'
    'D21' = '0.1313'
    'E21' = '2.01%'
    'D22' = '0.2311'
    'E22' = '-10.73%'
    'D23' = '0.04443'
    'E23' = '0.77%'
    'D24' = '0.001221'
    'E24' = '-1.23%'
    'D25' = '0.004361'
    'E25' = '-5.84%'
    'E26' = '-5.22%'
    'D27' = '0.0003990'
    'E27' = '0.03%'
    'D39' = '0.02510'
    'E39' = '11.42%'
    'D40' = '0.05258'
    'E40' = '3.01%'
    'D41' = '0.007714'
    'E41' = '3.25%'
    'D42' = '0.1431'
    'E42' = '5.75%'
    'D43' = '0.008607'
    'E43' = '-4.32%'
    'D44' = '0.002160'
    'E44' = '1.33%'
    'D45' = '0.008136'
    'E45' = '-5.69%'
    'D46' = '0.00006669'
    'E46' = '0.55%'
    'E47' = '0.03%'
    'B48' = 'CoinbaseStockToken'
    'C48' = 'https://coinranking.com/coin/_ZA6fIr53+coinbasestocktoken-coin'
    'D48' = '0.002483'
    'E48' = '148.31%'
    'B49' = 'BOLO'
    'C49' = 'https://coinranking.com/coin/ogrGe0dEab+bolo-bolo'
    'D49' = '0.002852'
    'E49' = '-0.29%'
    'D50' = '0.00002100'
    'E50' = '0.03%'
    'D51' = '0.0002000'
    'E51' = '0.03%'
}

foreach ($addr in $updates.Keys) {
    $value = $updates[$addr]
    $cell = $ws.Range($addr)
    if ($value -match '^-?[0-9.]+%?$') {
        # Numeric-looking text (prices / percentages): a leading apostrophe keeps
        # Excel from reinterpreting it as a number, which would normalise away
        # significant trailing zeros and the percent sign.
        $cell.Value = "'" + $value
        # The quote-prefix entry picks up an implicit number format; strip it so
        # the cell keeps the same (default) style it had before the edit.
        $cell.ClearFormats()
    } else {
        $cell.Value = $value
    }
}
